$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add guesstimate hours
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 4

# Update selection to B3, matching the author's last-selected cell
$ws.Range("B3").Select()
